# Atualizando site, criando novas funcionalidades
# Adds 5 new address entries (rows 27-31) to the collection sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Novo leblon ", "Rua"),
    @("Apraioh ", "Rua"),
    @("JH BEACH ", "Rua"),
    @("Ninah", "R"),
    @("TZ CORONEL", "Av José Luiz  Ferraz, 400 Bloco 1 - 306")
)

$startRow = 27
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
